# Auto-generated Excel COM-interop edit script
# Applies the cell-level numeric updates described in the commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 1629.3  # H17: was 1688.9143
$ws.Cells.Item(17, 9).Value = 900  # I17: was 0
$ws.Cells.Item(17, 10).Value = 1681.3928  # J17: was 1688.9143
$ws.Cells.Item(17, 11).Value = 2700  # K17: was 0
$ws.Cells.Item(17, 12).Value = 5044.178400000001  # L17: was 5066.742899999999
$ws.Cells.Item(17, 13).Value = -2532  # M17: was __ABSENT__
$ws.Cells.Item(17, 14).Value = -5380.178400000001  # N17: was -5402.742899999999
$ws.Cells.Item(138, 8).Value = 3091.86  # H138: was 3077.3232
$ws.Cells.Item(138, 9).Value = 1181.0667  # I138: was 1249.6207
$ws.Cells.Item(138, 10).Value = 3910.7715  # J138: was 3834.5144
$ws.Cells.Item(138, 11).Value = 3543.2001  # K138: was 3748.8621
$ws.Cells.Item(138, 12).Value = 11732.3145  # L138: was 11503.5432
$ws.Cells.Item(138, 13).Value = 1596.7999  # M138: was 1391.1379
$ws.Cells.Item(138, 14).Value = -22012.3145  # N138: was -21783.5432

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 43623.457  # H32: was 40810.95
$ws.Cells.Item(32, 9).Value = 42182.875  # I32: was 39252.77
$ws.Cells.Item(32, 11).Value = 42182.875  # K32: was 39252.77
$ws.Cells.Item(32, 13).Value = -41895.875  # M32: was -38965.77
$ws.Cells.Item(45, 8).Value = 1135.1818  # H45: was 1132
$ws.Cells.Item(45, 9).Value = 1061  # I45: was 1037.8
$ws.Cells.Item(45, 10).Value = 1333  # J45: was 1249.75
$ws.Cells.Item(45, 11).Value = 1061  # K45: was 1037.8
$ws.Cells.Item(45, 12).Value = 1333  # L45: was 1249.75
$ws.Cells.Item(45, 13).Value = -684  # M45: was -660.8
$ws.Cells.Item(45, 14).Value = -2087  # N45: was -2003.75
$ws.Cells.Item(61, 8).Value = 1846.8966  # H61: was 2179.8262
$ws.Cells.Item(61, 9).Value = 1713.8462  # I61: was 1979.7778
$ws.Cells.Item(61, 10).Value = 3000  # J61: was 2900
$ws.Cells.Item(61, 11).Value = 1713.8462  # K61: was 1979.7778
$ws.Cells.Item(61, 12).Value = 3000  # L61: was 2900
$ws.Cells.Item(61, 13).Value = -1501.8462  # M61: was -1767.7778
$ws.Cells.Item(61, 14).Value = -3424  # N61: was -3324
$ws.Cells.Item(74, 8).Value = 1780.421  # H74: was 1460.1864
$ws.Cells.Item(74, 9).Value = 1692.2667  # I74: was 1400.3864
$ws.Cells.Item(74, 10).Value = 2111  # J74: was 1635.6
$ws.Cells.Item(74, 11).Value = 1692.2667  # K74: was 1400.3864
$ws.Cells.Item(74, 12).Value = 2111  # L74: was 1635.6
$ws.Cells.Item(74, 13).Value = -818.2666999999999  # M74: was -526.3864000000001
$ws.Cells.Item(74, 14).Value = -3859  # N74: was -3383.6
$ws.Cells.Item(77, 8).Value = 1780.421  # H77: was 1460.1864
$ws.Cells.Item(77, 9).Value = 1692.2667  # I77: was 1400.3864
$ws.Cells.Item(77, 10).Value = 2111  # J77: was 1635.6
$ws.Cells.Item(77, 11).Value = 8461.333499999999  # K77: was 7001.932000000001
$ws.Cells.Item(77, 12).Value = 10555  # L77: was 8178
$ws.Cells.Item(77, 13).Value = -4093.333499999999  # M77: was -2633.932000000001
$ws.Cells.Item(77, 14).Value = -19291  # N77: was -16914
$ws.Cells.Item(115, 8).Value = 0  # H115: was 58000
$ws.Cells.Item(115, 10).Value = 0  # J115: was 58000
$ws.Cells.Item(115, 12).Value = 0  # L115: was 58000
$ws.Cells.Item(115, 14).ClearContents()  # N115: was -61134
$ws.Cells.Item(122, 8).Value = 1180.6  # H122: was 1238.2307
$ws.Cells.Item(122, 9).Value = 928.0909  # I122: was 955.2222
$ws.Cells.Item(122, 11).Value = 2784.2727  # K122: was 2865.6666
$ws.Cells.Item(122, 13).Value = -334.2727  # M122: was -415.6666
$ws.Cells.Item(132, 8).Value = 5848.456  # H132: was 6041.8545
$ws.Cells.Item(132, 9).Value = 7387.811  # I132: was 7779.6855
$ws.Cells.Item(132, 11).Value = 22163.433  # K132: was 23339.0565
$ws.Cells.Item(132, 13).Value = -19633.433  # M132: was -20809.0565
$ws.Cells.Item(136, 8).Value = 1846.8966  # H136: was 2179.8262
$ws.Cells.Item(136, 9).Value = 1713.8462  # I136: was 1979.7778
$ws.Cells.Item(136, 10).Value = 3000  # J136: was 2900
$ws.Cells.Item(136, 11).Value = 5141.5386  # K136: was 5939.3334
$ws.Cells.Item(136, 12).Value = 9000  # L136: was 8700
$ws.Cells.Item(136, 13).Value = -2591.5386  # M136: was -3389.3334
$ws.Cells.Item(136, 14).Value = -14100  # N136: was -13800

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 2416.4375  # H105: was 2788.5715
$ws.Cells.Item(105, 9).Value = 2237.6667  # I105: was 2630.9092
$ws.Cells.Item(105, 10).Value = 2952.75  # J105: was 3366.6667
$ws.Cells.Item(105, 11).Value = 2237.6667  # K105: was 2630.9092
$ws.Cells.Item(105, 12).Value = 2952.75  # L105: was 3366.6667
$ws.Cells.Item(105, 13).Value = -490.6667000000002  # M105: was -883.9092000000001
$ws.Cells.Item(105, 14).Value = -6446.75  # N105: was -6860.6667
$ws.Cells.Item(118, 8).Value = 0  # H118: was 42744.5
$ws.Cells.Item(118, 10).Value = 0  # J118: was 42744.5
$ws.Cells.Item(118, 12).Value = 0  # L118: was 42744.5
$ws.Cells.Item(118, 14).ClearContents()  # N118: was -46058.5
$ws.Cells.Item(134, 8).Value = 5543.073  # H134: was 6595.636
$ws.Cells.Item(134, 9).Value = 6041.1934  # I134: was 7724.6523
$ws.Cells.Item(134, 11).Value = 18123.5802  # K134: was 23173.9569
$ws.Cells.Item(134, 13).Value = -15588.5802  # M134: was -20638.9569

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 2036.6428  # H58: was 2071.1538
$ws.Cells.Item(58, 9).Value = 1875.5  # I58: was 1905.9131
$ws.Cells.Item(58, 10).Value = 3003.5  # J58: was 3338
$ws.Cells.Item(58, 11).Value = 1875.5  # K58: was 1905.9131
$ws.Cells.Item(58, 12).Value = 3003.5  # L58: was 3338
$ws.Cells.Item(58, 13).Value = -1672.5  # M58: was -1702.9131
$ws.Cells.Item(58, 14).Value = -3409.5  # N58: was -3744
$ws.Cells.Item(117, 8).Value = 48400  # H117: was 32857.145
$ws.Cells.Item(117, 10).Value = 48400  # J117: was 32857.145
$ws.Cells.Item(117, 12).Value = 48400  # L117: was 32857.145
$ws.Cells.Item(117, 14).Value = -57578  # N117: was -42035.145
$ws.Cells.Item(132, 8).Value = 4312213.5  # H132: was 2552339.5
$ws.Cells.Item(132, 9).Value = 1594  # I132: was 1041.9231
$ws.Cells.Item(132, 10).Value = 11365954  # J132: was 12502400
$ws.Cells.Item(132, 11).Value = 4782  # K132: was 3125.7693
$ws.Cells.Item(132, 12).Value = 34097862  # L132: was 37507200
$ws.Cells.Item(132, 13).Value = -2252  # M132: was -595.7692999999999
$ws.Cells.Item(132, 14).Value = -34102922  # N132: was -37512260
$ws.Cells.Item(134, 8).Value = 3861.6943  # H134: was 2470.61
$ws.Cells.Item(134, 9).Value = 4916.84  # I134: was 2585.608
$ws.Cells.Item(134, 10).Value = 1463.6364  # J134: was 1737.5
$ws.Cells.Item(134, 11).Value = 14750.52  # K134: was 7756.824000000001
$ws.Cells.Item(134, 12).Value = 4390.9092  # L134: was 5212.5
$ws.Cells.Item(134, 13).Value = -12215.52  # M134: was -5221.824000000001
$ws.Cells.Item(134, 14).Value = -9460.9092  # N134: was -10282.5
$ws.Cells.Item(136, 8).Value = 2036.6428  # H136: was 2071.1538
$ws.Cells.Item(136, 9).Value = 1875.5  # I136: was 1905.9131
$ws.Cells.Item(136, 10).Value = 3003.5  # J136: was 3338
$ws.Cells.Item(136, 11).Value = 5626.5  # K136: was 5717.7393
$ws.Cells.Item(136, 12).Value = 9010.5  # L136: was 10014
$ws.Cells.Item(136, 13).Value = -3076.5  # M136: was -3167.7393
$ws.Cells.Item(136, 14).Value = -14110.5  # N136: was -15114

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(35, 8).Value = 5500  # H35: was 3000
$ws.Cells.Item(35, 9).Value = 0  # I35: was 3000
$ws.Cells.Item(35, 10).Value = 5500  # J35: was 0
$ws.Cells.Item(35, 11).Value = 0  # K35: was 3000
$ws.Cells.Item(35, 12).Value = 5500  # L35: was 0
$ws.Cells.Item(35, 13).ClearContents()  # M35: was -2702
$ws.Cells.Item(35, 14).Value = -6096  # N35: was __ABSENT__
$ws.Cells.Item(102, 8).Value = 2232.2666  # H102: was 1874.2
$ws.Cells.Item(102, 9).Value = 2885.125  # I102: was 2083.1538
$ws.Cells.Item(102, 11).Value = 2885.125  # K102: was 2083.1538
$ws.Cells.Item(102, 13).Value = -1263.125  # M102: was -461.1538
$ws.Cells.Item(132, 8).Value = 5508.8  # H132: was 6016.0938
$ws.Cells.Item(132, 9).Value = 6156.84  # I132: was 6783.826
$ws.Cells.Item(132, 10).Value = 3888.7  # J132: was 4054.111
$ws.Cells.Item(132, 11).Value = 18470.52  # K132: was 20351.478
$ws.Cells.Item(132, 12).Value = 11666.1  # L132: was 12162.333
$ws.Cells.Item(132, 13).Value = -15940.52  # M132: was -17821.478
$ws.Cells.Item(132, 14).Value = -16726.1  # N132: was -17222.333

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 50001530  # H7: was 50002400
$ws.Cells.Item(7, 9).Value = 1784.8572  # I7: was 2000
$ws.Cells.Item(7, 10).Value = 166667600  # J7: was 71431144
$ws.Cells.Item(7, 11).Value = 1784.8572  # K7: was 2000
$ws.Cells.Item(7, 12).Value = 166667600  # L7: was 71431144
$ws.Cells.Item(7, 13).Value = -1672.8572  # M7: was -1888
$ws.Cells.Item(7, 14).Value = -166667824  # N7: was -71431368
$ws.Cells.Item(40, 8).Value = 3141.8125  # H40: was 3349.0833
$ws.Cells.Item(40, 9).Value = 2828.3845  # I40: was 2965.4443
$ws.Cells.Item(40, 11).Value = 2828.3845  # K40: was 2965.4443
$ws.Cells.Item(40, 13).Value = -2692.3845  # M40: was -2829.4443
$ws.Cells.Item(68, 8).Value = 1208  # H68: was 1280
$ws.Cells.Item(68, 10).Value = 1190  # J68: was 0
$ws.Cells.Item(68, 12).Value = 1190  # L68: was 0
$ws.Cells.Item(68, 14).Value = -2688  # N68: was __ABSENT__
$ws.Cells.Item(71, 8).Value = 1208  # H71: was 1280
$ws.Cells.Item(71, 10).Value = 1190  # J71: was 0
$ws.Cells.Item(71, 12).Value = 5950  # L71: was 0
$ws.Cells.Item(71, 14).Value = -13438  # N71: was __ABSENT__
$ws.Cells.Item(82, 8).Value = 2008.5834  # H82: was 1891.6428
$ws.Cells.Item(82, 10).Value = 2100.5  # J82: was 1872.875
$ws.Cells.Item(82, 12).Value = 2100.5  # L82: was 1872.875
$ws.Cells.Item(82, 14).Value = -2822.5  # N82: was -2594.875
$ws.Cells.Item(85, 8).Value = 2008.5834  # H85: was 1891.6428
$ws.Cells.Item(85, 10).Value = 2100.5  # J85: was 1872.875
$ws.Cells.Item(85, 12).Value = 2100.5  # L85: was 1872.875
$ws.Cells.Item(85, 14).Value = -4596.5  # N85: was -4368.875
$ws.Cells.Item(122, 8).Value = 2436.75  # H122: was 2800
$ws.Cells.Item(122, 9).Value = 1898.5  # I122: was 2000
$ws.Cells.Item(122, 10).Value = 2975  # J122: was 2960
$ws.Cells.Item(122, 11).Value = 5695.5  # K122: was 6000
$ws.Cells.Item(122, 12).Value = 8925  # L122: was 8880
$ws.Cells.Item(122, 13).Value = -3245.5  # M122: was -3550
$ws.Cells.Item(122, 14).Value = -13825  # N122: was -13780
$ws.Cells.Item(126, 8).Value = 50001530  # H126: was 50002400
$ws.Cells.Item(126, 9).Value = 1784.8572  # I126: was 2000
$ws.Cells.Item(126, 10).Value = 166667600  # J126: was 71431144
$ws.Cells.Item(126, 11).Value = 5354.571599999999  # K126: was 6000
$ws.Cells.Item(126, 12).Value = 500002800  # L126: was 214293432
$ws.Cells.Item(126, 13).Value = -2884.571599999999  # M126: was -3530
$ws.Cells.Item(126, 14).Value = -500007740  # N126: was -214298372
$ws.Cells.Item(132, 8).Value = 6760.4883  # H132: was 6767.4653
$ws.Cells.Item(132, 9).Value = 7325.0835  # I132: was 7505.8
$ws.Cells.Item(132, 10).Value = 3856.8572  # J132: was 3537.25
$ws.Cells.Item(132, 11).Value = 21975.2505  # K132: was 22517.4
$ws.Cells.Item(132, 12).Value = 11570.5716  # L132: was 10611.75
$ws.Cells.Item(132, 13).Value = -19445.2505  # M132: was -19987.4
$ws.Cells.Item(132, 14).Value = -16630.5716  # N132: was -15671.75
$ws.Cells.Item(136, 8).Value = 15445.333  # H136: was 7369.1577
$ws.Cells.Item(136, 9).Value = 18286.857  # I136: was 7942
$ws.Cells.Item(136, 10).Value = 5500  # J136: was 2500
$ws.Cells.Item(136, 11).Value = 54860.571  # K136: was 23826
$ws.Cells.Item(136, 12).Value = 16500  # L136: was 7500
$ws.Cells.Item(136, 13).Value = -52310.571  # M136: was -21276
$ws.Cells.Item(136, 14).Value = -21600  # N136: was -12600

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(109, 8).Value = 14947.368  # H109: was 15000
$ws.Cells.Item(109, 10).Value = 14947.368  # J109: was 15000
$ws.Cells.Item(109, 12).Value = 14947.368  # L109: was 15000
$ws.Cells.Item(109, 14).Value = -17721.368  # N109: was -17774
$ws.Cells.Item(132, 8).Value = 2272.0715  # H132: was 3564.4546
$ws.Cells.Item(132, 9).Value = 1960.4762  # I132: was 2958.7144
$ws.Cells.Item(132, 10).Value = 3206.8572  # J132: was 4624.5
$ws.Cells.Item(132, 11).Value = 5881.4286  # K132: was 8876.143199999999
$ws.Cells.Item(132, 12).Value = 9620.571599999999  # L132: was 13873.5
$ws.Cells.Item(132, 13).Value = -3351.4286  # M132: was -6346.143199999999
$ws.Cells.Item(132, 14).Value = -14680.5716  # N132: was -18933.5
